$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 27 (shifts old rows 27-31 down to 28-32)
$ws.Rows("27").Insert()

# 2. Copy formatting (styles only) from row 26 into the newly inserted row 27
$ws.Range("A26:E26").Copy()
$ws.Range("A27:E27").PasteSpecial(-4122)

# 3. Fix up row heights
$ws.Rows("23").RowHeight = 35.25
$ws.Rows("27").RowHeight = 22.5

# 4. Populate new cell text values (order chosen to match shared-string table append order)
$ws.Range("A20").Value = "گزارش های سیستم مالی"
$ws.Range("B22").Value = "تهیه اطلاعات گزارش دفتر روزنامه (مطابق ردیف های سند) در سرویس"
$ws.Range("B23").Value = "مشاهده اطلاعات گزارش دفتر روزنامه در فرم گزارشی جدید در برنامه"
$ws.Range("B25").Value = "تهیه اطلاعات گزارش دفتر حساب (مطابق ردیف های سند، برای مولفه حساب) در سرویس"
$ws.Range("B26").Value = "مشاهده اطلاعات گزارش دفتر حساب در فرم گزارشی جدید در برنامه"
$ws.Range("B24").Value = "طراحی گزارش دفتر روزنامه در محیط طراحی و یکپارچه سازی در برنامه"
$ws.Range("B27").Value = "طراحی گزارش دفتر حساب در محیط طراحی و یکپارچه سازی در برنامه"
$ws.Range("A28").Value = "رفع اشکالات موجود"
$ws.Range("B20").Value = "رفع اشکال گروه بندی حسابها در گزارش سند حسابداری - فرم مرسوم"
$ws.Range("B21").Value = "رفع اشکال گروه بندی حسابها در گزارش سند حسابداری - با سطوح شناور"

# 5. C column values for the new row (27) set to 1 (to match pattern of other data rows)
$ws.Range("C27").Value = 1

# 6. Resize the table (ListObject) to cover the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:E32"))

# 7. Fix data validation sqref: drop C26:C27 from the existing validated range so it becomes C28:C32
$ws.Range("C26:C27").Validation.Delete()

# 8. Update selection to match
$ws.Range("A27").Select()
